$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Faire s'affronter des unités :" paragraph gains a trailing sentence.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Faire s’affronter des unités :", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Faire s’affronter des unités : Créer des fonctions d’attaque qui sont appelées par les fonctions de portée.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert the three new "cahier des charges" paragraphs right after it,
#    before the blank paragraph + "Le Makefile" heading.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Faire s’affronter des unités :", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pIndex = $anchor.Paragraphs.Item(1).Index

$anchorRange = $d.Paragraphs.Item($pIndex).Range
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item($pIndex + 1)
$p1.Range.Text = "Finir le jeu et désigner un vainqueur : Vérifier à chaque tour de jeu si un des joueurs est mort et si le dernier tour est atteint. Si l’un des deux est mort, désigné l’autre vainqueur et si le nombre de tours max est atteint annoncer l’égalité."

$p2 = $d.Paragraphs.Item($pIndex + 2)
$p2.Range.Text = "Permettre de sauvegarder (Permettre de modifier facilement la sauvegarde) : Créer un fichier .txt ou supprimer celui qui porte déjà ce nom et écrire dedans ce qui est nécessaire à son fonctionnement. (Afficher en claire les informations pour que l’utilisateur puisse les modifier)"

# Paragraph 3 gets a temporary trailing char 'Z' so we can plant the _GoBack
# bookmark exactly after "Solution" (not at the paragraph-end boundary),
# then strip the helper char back out.
$p3 = $d.Paragraphs.Item($pIndex + 3)
$p3.Range.Text = "Objectif personnel : SolutionZ"
$p3start = $p3.Range.Start
$bmPos = $p3start + 29
$bmRange = $d.Range($bmPos, $bmPos)

# Remove the old "_GoBack" bookmark (currently sitting in the Makefile
# paragraph) before re-creating it at its new location.
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$zRange = $d.Range($bmPos, $bmPos + 1)
$zRange.Delete()

Write-Output "stage2-done"
